$d = $word.ActiveDocument

# --- 1. "Edin Mujagic*" -- remove spell-check proofing marks around "Mujagic" ---
$r = $d.Content
$r.Find.Execute("din Mujagic*,", $true, $false, $false, $false, $false, $true, 1, $false, "din Mujagic*,", 2) | Out-Null

# --- 2. "Alinea 3: (subthema)" -- remove spell-check proofing marks around "subthema" ---
$r = $d.Content
$r.Find.Execute("(subthema)", $true, $false, $false, $false, $false, $true, 1, $false, "(subthema)", 2) | Out-Null

# --- 3. "Kernzin alinea 3" paragraph -- rewrite + unbold + add quotes/semicolon/"en conglomeraten" ---
$r = $d.Content
$r.Find.Execute("Nederlanders voelen dit vooral in hun eigen portemonnees desondanks maken supermarkten hoge winsten. ", `
  $true, $false, $false, $false, $false, $true, 1, $false, `
  [char]34 + "Nederlanders voelen dit vooral in hun eigen portemonnees; desondanks maken supermarkten en conglomeraten hoge winsten." + [char]34, 2) | Out-Null

# Remove the bold formatting from that paragraph (content + paragraph mark)
$p = $d.Paragraphs(17)
$p.Range.Font.Bold = 0
$p.Range.Font.Bold = $false

# --- 4. "Resultaat van beleid ECB ... toekomst eruit zien?" -- remove proofing marks (graaiflatie / krimpflatie) ---
$r = $d.Content
$txt4 = "Resultaat van beleid ECB en krimp- & graaiflatie, deze punten hebben geleid tot de huidige situatie voor het afgelopen jaar. De prijzen van consumentengoederen daalt maar niet en er treden nu ook nieuwe economische fenomenen op zoals " + [char]0x2018 + "graaiflatie" + [char]0x2019 + " en " + [char]0x2018 + "krimpflatie" + [char]0x2019 + ". Hoe kunnen we dit herkennen en hoe gaat dit in de toekomst eruit zien?"
$r.Find.Execute($txt4, $true, $false, $false, $false, $false, $true, 1, $false, $txt4, 2) | Out-Null

# --- 5. "*bron Marco Kesteloo* " -- remove proofing marks around "Kesteloo" ---
$r = $d.Content
$txt5 = "*bron Marco Kesteloo* "
$r.Find.Execute($txt5, $true, $false, $false, $false, $false, $true, 1, $false, $txt5, 2) | Out-Null

# --- 6. "Supermarktketens weigeren jammergenoeg ... naar een " -- remove proofing marks around "jammergenoeg" ---
$r = $d.Content
$txt6 = "Supermarktketens weigeren jammergenoeg mijn vragen over verandering qua prijs in de supermarkt, maar in reactie naar een "
$r.Find.Execute($txt6, $true, $false, $false, $false, $false, $true, 1, $false, $txt6, 2) | Out-Null

# --- 7. "Initiatieven zoals foodwatch ... inhoudsafnames" -- remove proofing marks around "foodwatch" ---
$r = $d.Content
$txt7 = "Initiatieven zoals foodwatch bestuderen deze producten en tonen zowel prijsstijgingen als inhoudsafnames"
$r.Find.Execute($txt7, $true, $false, $false, $false, $false, $true, 1, $false, $txt7, 2) | Out-Null
